# Reorder the player roster rows (A2:C19) to match the updated roster order.
# Column A = Player, B = Position, C = Team. The header row (row 1) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Taurean Prince", "SG,SF", "Milwaukee Bucks"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Jalen Johnson", "PF", "Atlanta Hawks"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Royce O'Neale", "SF,PF", "Phoenix Suns"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Jalen Duren", "C", "Detroit Pistons")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value2 = $rec[0]
    $ws.Cells.Item($row, 2).Value2 = $rec[1]
    $ws.Cells.Item($row, 3).Value2 = $rec[2]
}
